# PICTURA_SQLPlan.xlsx -- "Added trigger TODOs" edit
#
# Adds 5 new trigger-related TODO rows to the "Triggers" sheet (rows 5-9,
# between the header row and the existing blank/query rows 11-14), marks
# each with a "TODO" status cell styled with the built-in "Bad" cell
# style, and leaves the workbook with the "Triggers" sheet as the active
# tab (selection resting on D9, the last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Triggers")

# Row 5 - single line
$ws.Range("C5").Value = "delete tag if no photo uses it anymore"
$ws.Range("E5").Value = "TODO"
$ws.Range("E5").Style = "Bad"
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("E5").VerticalAlignment = -4108

# Row 6 - wraps to two lines
$ws.Range("C6").Value = "comment only responds to comments on same photo"
$ws.Range("E6").Value = "TODO"
$ws.Range("E6").Style = "Bad"
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E6").VerticalAlignment = -4108
$ws.Rows.Item(6).RowHeight = 30

# Row 7 - single line
$ws.Range("C7").Value = "response's datetime > parent comment"
$ws.Range("E7").Value = "TODO"
$ws.Range("E7").Style = "Bad"
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").VerticalAlignment = -4108

# Row 8 - single line
$ws.Range("C8").Value = "comment datetime > photo datetime"
$ws.Range("E8").Value = "TODO"
$ws.Range("E8").Style = "Bad"
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E8").VerticalAlignment = -4108

# Row 9 - wraps to two lines
$ws.Range("C9").Value = "there must always be one admin per community"
$ws.Range("E9").Value = "TODO"
$ws.Range("E9").Style = "Bad"
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Rows.Item(9).RowHeight = 30

# Leave "Triggers" as the active sheet/tab, with the last cell touched
# (D9) selected -- matches the saved view state in the workbook.
$ws.Activate()
$ws.Range("D9").Select() | Out-Null
